# Applies the "updated DG and diagrams" edit: renames several class-diagram
# shapes from the AddressBook/Person model to the GradTrak/ModuleTaken model
# and nudges a handful of shape/connector geometries that shifted as a
# consequence of the longer label text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU = 12700  # points -> EMU

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# Converts a target EMU value to the points value to feed the COM layer,
# nudged by half an EMU so float32 truncation inside the host can't floor
# it down to the EMU below.
function Emu-ToPoints($emu) {
    return ($emu + 0.5) / $EMU
}

# "VersionedAddressBook" -> "VersionedGradTrak"
$shp46 = Get-ShapeById $s 46
$shp46.TextFrame.TextRange.Text = "VersionedGradTrak"

# "UniquePersonList" -> "UniqueModuleTakenList" (and box grows taller)
$shp49 = Get-ShapeById $s 49
$shp49.TextFrame.TextRange.Text = "UniqueModuleTakenList"
$shp49.Height = Emu-ToPoints 380656

# Connector between the two top-level boxes stretches with the resize
$shp30 = Get-ShapeById $s 30
$shp30.Height = Emu-ToPoints 18871
$shp30.Adjustments.Item(1) = 50000

# "Person" -> "ModuleTaken"
$shp62 = Get-ShapeById $s 62
$shp62.TextFrame.TextRange.Text = "ModuleTaken"

# "Name" -> "Code"
$shp76 = Get-ShapeById $s 76
$shp76.TextFrame.TextRange.Text = "Code"

# "Phone" -> "Semester"
$shp80 = Get-ShapeById $s 80
$shp80.TextFrame.TextRange.Text = "Semester"

# "Email" -> "Grade Range"
$shp83 = Get-ShapeById $s 83
$shp83.TextFrame.TextRange.Text = "Grade Range"

# "Address" -> "Workload" (box shifts left by 1 EMU and widens)
$shp85 = Get-ShapeById $s 85
$shp85.TextFrame.TextRange.Text = "Workload"
$shp85.Left = Emu-ToPoints 6680902
$shp85.Width = Emu-ToPoints 786697

# Connector feeding the "Workload" box widens by 1 EMU
$shp86 = Get-ShapeById $s 86
$shp86.Width = Emu-ToPoints 434401
$shp86.Adjustments.Item(1) = 50000
